$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("B2").Value = 1.097
$ws.Range("C2").Value = 95.4
$ws.Range("D2").Value = 0.63
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 0.41
$ws.Range("J2").Value = 632
$ws.Range("L2").Value = 394
$ws.Range("M2").Value = 1174
$ws.Range("O2").Value = 147
$ws.Range("Q2").Value = 1170
$ws.Range("S2").Value = 2.15
$ws.Range("T2").Value = -0.77
$ws.Range("U2").Value = -2.46
$ws.Range("V2").Value = 1.42
$ws.Range("W2").Value = 0.94
$ws.Range("X2").Value = 0.52
$ws.Range("Y2").Value = 0.81
$ws.Range("Z2").Value = 0.76

# --- Row 4 updates ---
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 0.41
$ws.Range("H4").Value = 0.51
$ws.Range("N4").Value = 624
$ws.Range("P4").Value = 443
$ws.Range("Q4").Value = 1212

# --- Row 5 updates ---
$ws.Range("D5").Value = 0.63
$ws.Range("G5").Value = 0.51
$ws.Range("N5").Value = 624
$ws.Range("Q5").Value = 1211

# --- Conditional formatting expression formula updates ---
$colThresholds = @{
    "J" = 624
    "K" = 149
    "L" = 391
    "M" = 1158
    "N" = 627
    "O" = 148
    "P" = 399
    "Q" = 1173
}

foreach ($col in $colThresholds.Keys) {
    $rng = $ws.Range("$col" + "2:" + "$col" + "5")
    $fcs = $rng.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fc = $fcs.Item($i)
        if ($fc.Type -eq 2) {
            $threshold = $colThresholds[$col]
            $fc.Formula1 = "=" + $col + "2<" + $threshold
        }
    }
}
